# Update the cryptos list values per the latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.886.59"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.639.11"
$ws.Range("E4").Value = "  -0.62%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "217.09"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("E8").Value = "  +1.68%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.0626"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +0.59%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "19.95"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +3.85%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.867.10"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.644.02"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.12"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -0.66%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.531"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +1.18%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "67.12"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").Value = "26.885.75"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  -0.23%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "219.79"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.51%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.91"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +3.80%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.41"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +0.88%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "2.45"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +3.81%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "9.18"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +0.17%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "146.93"
$cell.ClearFormats()
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("E26").Value = "  -0.62%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "7.36"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +3.47%  "
$ws.Range("E28").Value = "  +0.63%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "15.80"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +0.40%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.0505"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("E31").Value = "  -0.90%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.34"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -1.41%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.01"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +0.93%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.57"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("D35").Value = "1.262.31"
$ws.Range("E35").Value = "  -0.23%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.44"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("E37").Value = "  +2.29%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.536"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +1.02%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.835"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "1.777.85"
$ws.Range("E43").Value = "  -0.13%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "61.90"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("E45").Value = "  -1.55%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "91.63"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("E49").Value = "  -0.28%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "7.64"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +1.26%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0965"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -0.11%  "
